# Add a new week of "Cebolla" price observations (1a/2a/3a cosecha) for
# Agricola del Norte S.A. de Arica by inserting 3 rows right above the
# existing 2021-06-24 block (old row 588), pushing all the following rows
# down by 3 and growing the sheet from 609 to 612 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows before row 588; rows 588-609 shift down to 591-612.
$ws.Rows("588:590").Insert()

# Common (unchanged) field values shared by the three new rows.
$mercadoId = 1
$mercado   = "Agrícola del Norte S.A. de Arica"
$region    = "Arica y Parinacota"
$codreg    = 15
$catId     = 100112004
$categoria = "Cebolla"
$variedad  = "Sin especificar"
$unidad    = "`$/malla 18 kilos"
$origen    = "Región de Arica y Parinacota"
$kgUnid    = 18
$clasif    = "Hortaliza"
$fecha     = 44568

# Row 588: 1a (cosecha)
$ws.Cells.Item(588, 1).Value  = $mercadoId
$ws.Cells.Item(588, 2).Value  = $mercado
$ws.Cells.Item(588, 3).Value  = $region
$ws.Cells.Item(588, 4).Value  = $fecha
$ws.Cells.Item(588, 5).Value  = $codreg
$ws.Cells.Item(588, 6).Value  = $catId
$ws.Cells.Item(588, 7).Value  = $categoria
$ws.Cells.Item(588, 8).Value  = $variedad
$ws.Cells.Item(588, 9).Value  = "1a (cosecha)"
$ws.Cells.Item(588, 10).Value = 300
$ws.Cells.Item(588, 11).Value = 2500
$ws.Cells.Item(588, 12).Value = 3000
$ws.Cells.Item(588, 13).Value = 2750
$ws.Cells.Item(588, 14).Value = $unidad
$ws.Cells.Item(588, 15).Value = $origen
$ws.Cells.Item(588, 16).Value = 153
$ws.Cells.Item(588, 17).Value = $kgUnid
$ws.Cells.Item(588, 18).Value = $clasif

# Row 589: 2a (cosecha)
$ws.Cells.Item(589, 1).Value  = $mercadoId
$ws.Cells.Item(589, 2).Value  = $mercado
$ws.Cells.Item(589, 3).Value  = $region
$ws.Cells.Item(589, 4).Value  = $fecha
$ws.Cells.Item(589, 5).Value  = $codreg
$ws.Cells.Item(589, 6).Value  = $catId
$ws.Cells.Item(589, 7).Value  = $categoria
$ws.Cells.Item(589, 8).Value  = $variedad
$ws.Cells.Item(589, 9).Value  = "2a (cosecha)"
$ws.Cells.Item(589, 10).Value = 350
$ws.Cells.Item(589, 11).Value = 2000
$ws.Cells.Item(589, 12).Value = 2500
$ws.Cells.Item(589, 13).Value = 2250
$ws.Cells.Item(589, 14).Value = $unidad
$ws.Cells.Item(589, 15).Value = $origen
$ws.Cells.Item(589, 16).Value = 125
$ws.Cells.Item(589, 17).Value = $kgUnid
$ws.Cells.Item(589, 18).Value = $clasif

# Row 590: 3a (cosecha)
$ws.Cells.Item(590, 1).Value  = $mercadoId
$ws.Cells.Item(590, 2).Value  = $mercado
$ws.Cells.Item(590, 3).Value  = $region
$ws.Cells.Item(590, 4).Value  = $fecha
$ws.Cells.Item(590, 5).Value  = $codreg
$ws.Cells.Item(590, 6).Value  = $catId
$ws.Cells.Item(590, 7).Value  = $categoria
$ws.Cells.Item(590, 8).Value  = $variedad
$ws.Cells.Item(590, 9).Value  = "3a (cosecha)"
$ws.Cells.Item(590, 10).Value = 300
$ws.Cells.Item(590, 11).Value = 1500
$ws.Cells.Item(590, 12).Value = 2000
$ws.Cells.Item(590, 13).Value = 1750
$ws.Cells.Item(590, 14).Value = $unidad
$ws.Cells.Item(590, 15).Value = $origen
$ws.Cells.Item(590, 16).Value = 97
$ws.Cells.Item(590, 17).Value = $kgUnid
$ws.Cells.Item(590, 18).Value = $clasif
